# Bonsai On A Budget — Art & Science slide:
# Split the opening sentence of the content placeholder into three runs so
# that the new phrase "living trompe-l'oeil" is inserted and shown in bold:
#   "Bonsai: a small tree that appears to be a scaled-down large tree."
#   -> "A bonsai is a " + "living trompe-l'oeil" (bold) + ": a small tree
#      that appears to be a scaled-down large tree."

$p = $ppt.ActivePresentation

$oldText = "Bonsai: a small tree that appears to be a scaled-down large tree."

# Locate the shape (on any slide) whose text contains the sentence we need
# to edit, rather than relying on a hard-coded slide/shape index.
$targetShape = $null
$targetSlide = $null

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText -and ($shp.TextFrame.TextRange.Text -like "*$oldText*")) {
                $targetShape = $shp
                $targetSlide = $slide
            }
        }
    }
}

if ($targetShape -eq $null) {
    Write-Output "ERROR: could not find shape containing target sentence"
} else {
    $tr = $targetShape.TextFrame.TextRange

    # Find which paragraph holds the sentence (it is the first one, but look
    # it up defensively instead of assuming).
    $paraCount = $tr.Paragraphs().Count
    $targetParaIndex = -1
    for ($k = 1; $k -le $paraCount; $k++) {
        $pk = $tr.Paragraphs($k, 1)
        if ($pk.Text -like "$oldText*") {
            $targetParaIndex = $k
        }
    }

    $para = $tr.Paragraphs($targetParaIndex, 1)

    $prefix = "A bonsai is a "
    $bolded = "living trompe-l" + [char]0x2019 + "oeil"
    $suffix = ": a small tree that appears to be a scaled-down large tree."

    # Setting .Text on the whole paragraph keeps it as a single run with the
    # paragraph's existing run-level formatting (matches the non-bold runs
    # either side of the new phrase).
    $para.Text = $prefix + $bolded + $suffix

    # Now bold just the newly-inserted phrase; re-applying a distinct
    # character format over part of the run range splits it into its own
    # <a:r> run, exactly like typing then selecting + bolding in the UI.
    $boldStart = $para.Start + $prefix.Length
    $boldLen = $bolded.Length
    $boldRange = $tr.Characters($boldStart, $boldLen)
    $boldRange.Font.Bold = 1
}
